$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 38
$ws.Range("I11").Value = 38
$ws.Range("K11").Value = 38
$ws.Range("M11").Value = 102
$ws.Range("H53").Value = 52632556
$ws.Range("I53").Value = 71429500
$ws.Range("K53").Value = 71429500
$ws.Range("M53").Value = -71428863
$ws.Range("H98").Value = 2999.5
$ws.Range("I98").Value = 2999.5
$ws.Range("K98").Value = 2999.5
$ws.Range("M98").Value = -1501.5
$ws.Range("H122").Value = 2999.5
$ws.Range("I122").Value = 2999.5
$ws.Range("K122").Value = 8998.5
$ws.Range("M122").Value = -6548.5
$ws.Range("H132").Value = 4845.6
$ws.Range("I132").Value = 5552.407
$ws.Range("K132").Value = 16657.221
$ws.Range("M132").Value = -14127.221
$ws.Range("H137").Value = 5002.8096
$ws.Range("I137").Value = 2579.6
$ws.Range("J137").Value = 7205.727
$ws.Range("K137").Value = 7738.799999999999
$ws.Range("L137").Value = 21617.181
$ws.Range("M137").Value = -5188.799999999999
$ws.Range("N137").Value = -26717.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2290695
$ws.Range("I32").Value = 1237638.8
$ws.Range("K32").Value = 1237638.8
$ws.Range("M32").Value = -1237351.8
$ws.Range("H61").Value = 5044.8184
$ws.Range("I61").Value = 4732.467
$ws.Range("K61").Value = 4732.467
$ws.Range("M61").Value = -4520.467
$ws.Range("H74").Value = 24527824
$ws.Range("I74").Value = 193241.89
$ws.Range("J74").Value = 83336390
$ws.Range("K74").Value = 193241.89
$ws.Range("L74").Value = 83336390
$ws.Range("M74").Value = -192367.89
$ws.Range("N74").Value = -83338138
$ws.Range("H77").Value = 24527824
$ws.Range("I77").Value = 193241.89
$ws.Range("J77").Value = 83336390
$ws.Range("K77").Value = 966209.4500000001
$ws.Range("L77").Value = 416681950
$ws.Range("M77").Value = -961841.4500000001
$ws.Range("N77").Value = -416690686
$ws.Range("H102").Value = 1714.2307
$ws.Range("I102").Value = 1507.7273
$ws.Range("K102").Value = 1507.7273
$ws.Range("M102").Value = 114.2727
$ws.Range("H110").Value = 1159.4814
$ws.Range("I110").Value = 1018.9545
$ws.Range("K110").Value = 1018.9545
$ws.Range("M110").Value = 1026.0455
$ws.Range("H136").Value = 5044.8184
$ws.Range("I136").Value = 4732.467
$ws.Range("K136").Value = 14197.401
$ws.Range("M136").Value = -11647.401

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1944.1428
$ws.Range("I134").Value = 1578.1177
$ws.Range("K134").Value = 4734.3531
$ws.Range("M134").Value = -2199.3531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 839.86664
$ws.Range("I7").Value = 987.4167
$ws.Range("K7").Value = 987.4167
$ws.Range("M7").Value = -874.4167
$ws.Range("H22").Value = 776
$ws.Range("J22").Value = 739.6667
$ws.Range("L22").Value = 739.6667
$ws.Range("N22").Value = -1439.6667
$ws.Range("H31").Value = 1788619.6
$ws.Range("I31").Value = 1269.1
$ws.Range("K31").Value = 1269.1
$ws.Range("M31").Value = -974.0999999999999
$ws.Range("H34").Value = 1788619.6
$ws.Range("I34").Value = 1269.1
$ws.Range("K34").Value = 1269.1
$ws.Range("M34").Value = -1067.1
$ws.Range("H36").Value = 25023.5
$ws.Range("I36").Value = 10048
$ws.Range("K36").Value = 10048
$ws.Range("M36").Value = -9660
$ws.Range("H40").Value = 25023.5
$ws.Range("I40").Value = 10048
$ws.Range("K40").Value = 10048
$ws.Range("M40").Value = -9888
$ws.Range("H58").Value = 4977.815
$ws.Range("I58").Value = 2368.875
$ws.Range("J58").Value = 8772.637000000001
$ws.Range("K58").Value = 2368.875
$ws.Range("L58").Value = 8772.637000000001
$ws.Range("M58").Value = -2165.875
$ws.Range("N58").Value = -9178.637000000001
$ws.Range("H132").Value = 3779.6553
$ws.Range("I132").Value = 3302.4707
$ws.Range("K132").Value = 9907.4121
$ws.Range("M132").Value = -7377.4121
$ws.Range("H134").Value = 4740.2856
$ws.Range("I134").Value = 4941.5
$ws.Range("J134").Value = 3533
$ws.Range("K134").Value = 14824.5
$ws.Range("L134").Value = 10599
$ws.Range("M134").Value = -12289.5
$ws.Range("N134").Value = -15669
$ws.Range("H136").Value = 4977.815
$ws.Range("I136").Value = 2368.875
$ws.Range("J136").Value = 8772.637000000001
$ws.Range("K136").Value = 7106.625
$ws.Range("L136").Value = 26317.911
$ws.Range("M136").Value = -4556.625
$ws.Range("N136").Value = -31417.911

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1099.7391
$ws.Range("I2").Value = 211
$ws.Range("J2").Value = 2482.2222
$ws.Range("K2").Value = 1266
$ws.Range("L2").Value = 14893.3332
$ws.Range("M2").Value = -1153
$ws.Range("N2").Value = -15119.3332
$ws.Range("H60").Value = 1054317
$ws.Range("I60").Value = 5000325
$ws.Range("J60").Value = 2048.3333
$ws.Range("K60").Value = 15000975
$ws.Range("L60").Value = 6144.999899999999
$ws.Range("M60").Value = -15000724
$ws.Range("N60").Value = -6646.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 55046
$ws.Range("J46").Value = 55046
$ws.Range("L46").Value = 55046
$ws.Range("N46").Value = -55358
$ws.Range("H70").Value = 71435290
$ws.Range("I70").Value = 125005000
$ws.Range("K70").Value = 125005000
$ws.Range("M70").Value = -125004730
$ws.Range("H73").Value = 71435290
$ws.Range("I73").Value = 125005000
$ws.Range("K73").Value = 125005000
$ws.Range("M73").Value = -125004064
$ws.Range("H102").Value = 37040860
$ws.Range("I102").Value = 40003812
$ws.Range("K102").Value = 40003812
$ws.Range("M102").Value = -40002190
$ws.Range("H132").Value = 2383.1853
$ws.Range("I132").Value = 1938.5
$ws.Range("K132").Value = 5815.5
$ws.Range("M132").Value = -3285.5
$ws.Range("H135").Value = 72407.336
$ws.Range("J135").Value = 72407.336
$ws.Range("L135").Value = 72407.336
$ws.Range("N135").Value = -82547.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 14941414
$ws.Range("J2").Value = 20823980
$ws.Range("L2").Value = 20823980
$ws.Range("N2").Value = -20824204
$ws.Range("H40").Value = 43942.2
$ws.Range("I40").Value = 69497.44500000001
$ws.Range("K40").Value = 69497.44500000001
$ws.Range("M40").Value = -69361.44500000001
$ws.Range("H55").Value = 304.27274
$ws.Range("I55").Value = 279.29413
$ws.Range("K55").Value = 279.29413
$ws.Range("M55").Value = -106.29413
$ws.Range("H114").Value = 119000
$ws.Range("J114").Value = 119000
$ws.Range("L114").Value = 119000
$ws.Range("N114").Value = -127678
$ws.Range("H120").Value = 200000
$ws.Range("J120").Value = 200000
$ws.Range("L120").Value = 200000
$ws.Range("N120").Value = -209676
$ws.Range("H122").Value = 4654.75
$ws.Range("I122").Value = 4659.1665
$ws.Range("J122").Value = 4650.3335
$ws.Range("K122").Value = 13977.4995
$ws.Range("L122").Value = 13951.0005
$ws.Range("M122").Value = -11527.4995
$ws.Range("N122").Value = -18851.0005
$ws.Range("H132").Value = 4358.1514
$ws.Range("I132").Value = 4675.727
$ws.Range("K132").Value = 14027.181
$ws.Range("M132").Value = -11497.181
$ws.Range("H136").Value = 6551.273
$ws.Range("I136").Value = 6066
$ws.Range("J136").Value = 6733.25
$ws.Range("K136").Value = 18198
$ws.Range("L136").Value = 20199.75
$ws.Range("M136").Value = -15648
$ws.Range("N136").Value = -25299.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8335800.5
$ws.Range("I122").Value = 2582.8635
$ws.Range("K122").Value = 7748.5905
$ws.Range("M122").Value = -5298.5905
$ws.Range("H132").Value = 1814.7872
$ws.Range("I132").Value = 1625.8422
$ws.Range("J132").Value = 2612.5557
$ws.Range("K132").Value = 4877.5266
$ws.Range("L132").Value = 7837.6671
$ws.Range("M132").Value = -2347.5266
$ws.Range("N132").Value = -12897.6671
$ws.Range("H136").Value = 11911794
$ws.Range("I136").Value = 15880000
$ws.Range("J136").Value = 7174.4287
$ws.Range("K136").Value = 47640000
$ws.Range("L136").Value = 21523.2861
$ws.Range("M136").Value = -47637450
$ws.Range("N136").Value = -26623.2861
